{"js": "// Adds a new \"mailing address\" placeholder line under the existing\n// {{propertyAddress}} line (near the top of the letter, outside the\n// table), turning it into {{propertyAddress_st_unit}} followed by a new\n// paragraph {{propertyAddress_city_state_zip}}. Also removes a stray\n// empty \"No Spacing\" paragraph that sat right after the\n// \"{{associationName}} Board of Directors\" line.\n\nconst body = context.document.body;\n\n// --- Part 1: split {{propertyAddress}} into the street/unit line and a\n// new city/state/zip line -----------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph whose *entire* text is exactly \"{{propertyAddress}}\".\n// (The same placeholder also appears once inside the table further down\n// the letter - that occurrence must stay untouched.)\nlet targetPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"{{propertyAddress}}\") {\n    targetPara = p;\n    break;\n  }\n}\n\nif (targetPara) {\n  // Append \"_st_unit\" right after the \"propertyAddress\" word so the line\n  // reads \"{{propertyAddress_st_unit}}\".\n  const hit = targetPara.search(\"propertyAddress\", { matchCase: true });\n  hit.load(\"items\");\n  await context.sync();\n  if (hit.items.length > 0) {\n    hit.items[0].insertText(\"_st_unit\", \"After\");\n    await context.sync();\n  }\n\n  // Insert a brand-new paragraph right after it for the city/state/zip\n  // placeholder; it naturally inherits the same paragraph formatting.\n  targetPara.insertParagraph(\"{{propertyAddress_city_state_zip}}\", \"After\");\n  await context.sync();\n}\n\n// --- Part 2: drop the empty \"No Spacing\" paragraph directly below\n// \"{{associationName}} Board of Directors\" --------------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const p = paragraphs2.items[i];\n  if (p.text.indexOf(\"Board of Directors\") !== -1) {\n    const next = paragraphs2.items[i + 1];\n    if (next && next.text === \"\" && next.style === \"No Spacing\") {\n      next.delete();\n    }\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Adds a new \"mailing address\" placeholder line under the existing\n# {{propertyAddress}} line (near the top of the letter, outside the\n# table), turning it into {{propertyAddress_st_unit}} followed by a new\n# paragraph {{propertyAddress_city_state_zip}}. Also removes a stray\n# empty \"No Spacing\" paragraph that sat right after the\n# \"{{associationName}} Board of Directors\" line.\n\n$d = $word.ActiveDocument\n\n# --- Part 1: split {{propertyAddress}} into the street/unit line and a\n# new city/state/zip line -----------------------------------------\n\n# Find the paragraph whose *entire* text is exactly \"{{propertyAddress}}\"\n# followed by the paragraph mark. (The same placeholder also appears once\n# inside the table further down the letter - that occurrence must stay\n# untouched, so we match on the body-level Paragraphs collection only,\n# stopping at the first hit.)\n$targetPara = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq (\"{{propertyAddress}}\" + [char]13)) {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    # Append \"_st_unit\" right after the \"propertyAddress\" word so the\n    # line reads \"{{propertyAddress_st_unit}}\".\n    $hit = $targetPara.Range\n    $find = $hit.Find\n    $find.Text = \"propertyAddress\"\n    $find.Forward = $true\n    $find.Wrap = 0\n    $found = $find.Execute()\n    if ($found) {\n        $hit.InsertAfter(\"_st_unit\")\n    }\n\n    # Insert a brand-new paragraph right after it for the city/state/zip\n    # placeholder; it naturally inherits the same paragraph formatting.\n    $endRng = $targetPara.Range\n    $endRng.Collapse(0)\n    $endRng.InsertParagraphAfter()\n    $endRng.InsertAfter(\"{{propertyAddress_city_state_zip}}\")\n}\n\n# --- Part 2: drop the empty \"No Spacing\" paragraph directly below\n# \"{{associationName}} Board of Directors\" --------------------------\n$count2 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Board of Directors*\") {\n        if ($i + 1 -le $count2) {\n            $next = $d.Paragraphs.Item($i + 1)\n            if ($next.Range.Text -eq [char]13 -and $next.Style.NameLocal -eq \"No Spacing\") {\n                $next.Range.Delete()\n            }\n        }\n        break\n    }\n}\n"}
